$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The backup-code values that used to live in A10:A12 move up to A2:A4
# (replacing the old A2:A4 codes), and A10:A12 are cleared out. The
# codes in A13:A15 are untouched.
$ws.Range("A2").Value = "C9DEXVFAR31A"
$ws.Range("A3").Value = "05ANDJ337D9B"
$ws.Range("A4").Value = "SAS5DZQK4GHR"

$ws.Range("A10:A12").ClearContents()

# Update the active selection to match the saved workbook state.
$ws.Range("A6").Select()
